$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 628
$ws.Range("F3").Value = 76
$ws.Range("F5").Value = 4585
$ws.Range("F6").Value = 1826
$ws.Range("F7").Value = 129
$ws.Range("F9").Value = 3084
$ws.Range("F11").Value = 583
$ws.Range("F12").Value = 246
$ws.Range("F13").Value = 604
$ws.Range("F14").Value = 519
$ws.Range("F15").Value = 517
$ws.Range("F16").Value = 368
$ws.Range("F18").Value = 1764
$ws.Range("F19").Value = 1315
$ws.Range("F21").Value = 1566
$ws.Range("F23").Value = 606
$ws.Range("F25").Value = 528
$ws.Range("F27").Value = 45
$ws.Range("F28").Value = 92
$ws.Range("F30").Value = 83
$ws.Range("F31").Value = 3641
$ws.Range("F32").Value = 752
$ws.Range("F34").Value = 295
$ws.Range("F35").Value = 54
$ws.Range("F36").Value = 1745

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 20

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 628
$ws.Range("F3").Value = 76
$ws.Range("F5").Value = 4585
$ws.Range("F6").Value = 1826
$ws.Range("F7").Value = 129
$ws.Range("F9").Value = 3084
$ws.Range("F11").Value = 583
$ws.Range("F12").Value = 246
$ws.Range("F13").Value = 604
$ws.Range("F14").Value = 519
$ws.Range("F15").Value = 517
$ws.Range("F16").Value = 20
$ws.Range("F17").Value = 368
$ws.Range("F19").Value = 1764
$ws.Range("F20").Value = 1315
$ws.Range("F22").Value = 1566
$ws.Range("F24").Value = 606
$ws.Range("F26").Value = 528
$ws.Range("F28").Value = 45
$ws.Range("F29").Value = 92
$ws.Range("F31").Value = 83
$ws.Range("F32").Value = 3641
$ws.Range("F34").Value = 752
$ws.Range("F36").Value = 295
$ws.Range("F37").Value = 54
$ws.Range("F38").Value = 1745
